$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows at row 210, shifting existing rows 210-278 down to 215-283
$ws.Range("A210:T214").EntireRow.Insert()

# Row 210
$ws.Cells.Item(210, 1).Value = 10
$ws.Cells.Item(210, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(210, 3).Value = 'La Araucanía'
$ws.Cells.Item(210, 4).Value = 44636
$ws.Cells.Item(210, 5).Value = 9
$ws.Cells.Item(210, 6).Value = 'Fruta'
$ws.Cells.Item(210, 7).Value = 100103
$ws.Cells.Item(210, 8).Value = 'Frutos de hueso (carozo)'
$ws.Cells.Item(210, 9).Value = 100103004
$ws.Cells.Item(210, 10).Value = 'Durazno'
$ws.Cells.Item(210, 11).Value = 'Phillips Cling'
$ws.Cells.Item(210, 12).Value = 'Primera'
$ws.Cells.Item(210, 13).Value = 210
$ws.Cells.Item(210, 14).Value = 16000
$ws.Cells.Item(210, 15).Value = 16000
$ws.Cells.Item(210, 16).Value = 16000
$ws.Cells.Item(210, 17).Value = '$/bandeja 18 kilos granel'
$ws.Cells.Item(210, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(210, 19).Value = 889
$ws.Cells.Item(210, 20).Value = 18

# Row 211
$ws.Cells.Item(211, 1).Value = 10
$ws.Cells.Item(211, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(211, 3).Value = 'La Araucanía'
$ws.Cells.Item(211, 4).Value = 44636
$ws.Cells.Item(211, 5).Value = 9
$ws.Cells.Item(211, 6).Value = 'Fruta'
$ws.Cells.Item(211, 7).Value = 100103
$ws.Cells.Item(211, 8).Value = 'Frutos de hueso (carozo)'
$ws.Cells.Item(211, 9).Value = 100103004
$ws.Cells.Item(211, 10).Value = 'Durazno'
$ws.Cells.Item(211, 11).Value = 'September Sweet'
$ws.Cells.Item(211, 12).Value = 'Especial'
$ws.Cells.Item(211, 13).Value = 210
$ws.Cells.Item(211, 14).Value = 20000
$ws.Cells.Item(211, 15).Value = 20000
$ws.Cells.Item(211, 16).Value = 20000
$ws.Cells.Item(211, 17).Value = '$/bandeja 18 kilos empedrada'
$ws.Cells.Item(211, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(211, 19).Value = 1111
$ws.Cells.Item(211, 20).Value = 18

# Row 212
$ws.Cells.Item(212, 1).Value = 10
$ws.Cells.Item(212, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(212, 3).Value = 'La Araucanía'
$ws.Cells.Item(212, 4).Value = 44636
$ws.Cells.Item(212, 5).Value = 9
$ws.Cells.Item(212, 6).Value = 'Fruta'
$ws.Cells.Item(212, 7).Value = 100103
$ws.Cells.Item(212, 8).Value = 'Frutos de hueso (carozo)'
$ws.Cells.Item(212, 9).Value = 100103004
$ws.Cells.Item(212, 10).Value = 'Durazno'
$ws.Cells.Item(212, 11).Value = 'September Sweet'
$ws.Cells.Item(212, 12).Value = 'Primera'
$ws.Cells.Item(212, 13).Value = 255
$ws.Cells.Item(212, 14).Value = 16000
$ws.Cells.Item(212, 15).Value = 16000
$ws.Cells.Item(212, 16).Value = 16000
$ws.Cells.Item(212, 17).Value = '$/bandeja 18 kilos granel'
$ws.Cells.Item(212, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(212, 19).Value = 889
$ws.Cells.Item(212, 20).Value = 18

# Row 213
$ws.Cells.Item(213, 1).Value = 10
$ws.Cells.Item(213, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(213, 3).Value = 'La Araucanía'
$ws.Cells.Item(213, 4).Value = 44636
$ws.Cells.Item(213, 5).Value = 9
$ws.Cells.Item(213, 6).Value = 'Fruta'
$ws.Cells.Item(213, 7).Value = 100103
$ws.Cells.Item(213, 8).Value = 'Frutos de hueso (carozo)'
$ws.Cells.Item(213, 9).Value = 100103004
$ws.Cells.Item(213, 10).Value = 'Durazno'
$ws.Cells.Item(213, 11).Value = 'September Sweet'
$ws.Cells.Item(213, 12).Value = 'Primera'
$ws.Cells.Item(213, 13).Value = 6
$ws.Cells.Item(213, 14).Value = 460000
$ws.Cells.Item(213, 15).Value = 460000
$ws.Cells.Item(213, 16).Value = 460000
$ws.Cells.Item(213, 17).Value = '$/bins (400 kilos)'
$ws.Cells.Item(213, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(213, 19).Value = 1150
$ws.Cells.Item(213, 20).Value = 400

# Row 214
$ws.Cells.Item(214, 1).Value = 10
$ws.Cells.Item(214, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(214, 3).Value = 'La Araucanía'
$ws.Cells.Item(214, 4).Value = 44636
$ws.Cells.Item(214, 5).Value = 9
$ws.Cells.Item(214, 6).Value = 'Fruta'
$ws.Cells.Item(214, 7).Value = 100103
$ws.Cells.Item(214, 8).Value = 'Frutos de hueso (carozo)'
$ws.Cells.Item(214, 9).Value = 100103004
$ws.Cells.Item(214, 10).Value = 'Durazno'
$ws.Cells.Item(214, 11).Value = 'September Sweet'
$ws.Cells.Item(214, 12).Value = 'Segunda'
$ws.Cells.Item(214, 13).Value = 210
$ws.Cells.Item(214, 14).Value = 14000
$ws.Cells.Item(214, 15).Value = 14000
$ws.Cells.Item(214, 16).Value = 14000
$ws.Cells.Item(214, 17).Value = '$/bandeja 18 kilos granel'
$ws.Cells.Item(214, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(214, 19).Value = 778
$ws.Cells.Item(214, 20).Value = 18
